$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "num"/"formulae" helper columns (F, G) were a scratch calc that's no
# longer needed - clear their header text, values and the SUM() formulas.
$ws.Range("F1:G3").ClearContents()

# A new Google-Form-style response row came in (auto user add) - replicate
# the formatting of the previous response row onto the new one, then fill
# in the values.
$ws.Range("A3:E3").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)

$ws.Range("A4").Value = 44967.7536921296
$ws.Range("B4").Value = "Lady"
$ws.Range("C4").Value = "Gaga"
$ws.Range("D4").Value = 2024
$ws.Range("E4").Value = "queen@wpi.edu"

# Turn the e-mail address into a live mailto: hyperlink, matching the other
# rows' WPI-email column.
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:queen@wpi.edu", "", "", "queen@wpi.edu")

# Adding the hyperlink re-stamps the cell with Excel's builtin "Hyperlink"
# style (underline + theme color) - put the plain blue-font style used by
# the sibling cells back, then drop the now-unused builtin style.
$ws.Range("E3").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$wb.Styles.Item("Hyperlink").Delete()

[void]$ws.Range("D31").Select()
